$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Problem or Target Behavior"
$ws.Range("C1").Value = "Score"

# Data
$dates = @("2024-05-22","2024-05-23","2024-05-24","2024-05-25","2024-05-26","2024-05-27","2024-05-28","2024-05-29")
$scores = @(5,8,7,7,6,6,7,5)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]

    $ws.Cells.Item($row, 2).Value = "Stress"
    $ws.Cells.Item($row, 3).Value = $scores[$i]
}

# Rows 7-9 are new; give column A the same date number format already used
# by the existing rows (2-6) so we reuse the existing style instead of
# creating a fresh custom number format.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("A7:A9").PasteSpecial(-4122) | Out-Null

# Underline the C2 score cell
$ws.Cells.Item(2, 3).Font.Underline = $true

# Column A now holds dates and needs its own best-fit width (it previously
# had no explicit width, relying on the sheet default).
$ws.Columns.Item(1).AutoFit() | Out-Null

# Selection matching final state
$ws.Range("C10").Select()
